$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1
$ws.Range("C4").Value = 1.1
$ws.Range("B5").Value = -0.1
$ws.Range("C7").Value = 1.1

$ws.Range("C8").Select()
